$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.6961823105812073
$ws.Range("B1").Value = 1.327176451683044
$ws.Range("C1").Value = 3.717926502227783
$ws.Range("D1").Value = 2.666195869445801
$ws.Range("E1").Value = 0.5914642214775085
